# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps to reflect the new report run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date (same value also shown on
# the de-de sheet as Correspond Handoff Datetime)
$overview.Range("G2").Value = "2016-09-03 05:11:20"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$zhcn.Range("H2").Value = "2016-09-03 05:11:15"
$zhcn.Range("K2").Value = "2016-09-03 05:11:32"

# de-de sheet: Correspond Handoff Datetime (mirrors Overview!G2) / Correspond Handback DateTime
$dede.Range("H2").Value = "2016-09-03 05:11:20"
$dede.Range("K2").Value = "2016-09-03 05:11:40"
